$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 text changes from "Usuário" to "Username"
$ws.Range("A1").Value = "Username"

# Remove the leftover styled-but-empty row (row 4 only held B4, an
# empty cell that just carried a font style) so the used range shrinks
# back down to A1:B2, matching the target sheet.
$ws.Rows(4).Delete()

# Set a custom width on column A (closest representable value to 9.85546875)
$ws.Columns("A").ColumnWidth = 9

# Move the active selection to D3
$null = $ws.Range("D3").Select()
